$wb = $excel.ActiveWorkbook

# --- StatOutput sheet: row 2 counts change (breed filter switched from
#     'Akita' to 'German Shepherd Dog' -> different result counts) ---
$statWs = $wb.Worksheets.Item("StatOutput")

# number_of_files: 1 -> 20
$statWs.Range("A2").NumberFormat = "@"
$statWs.Range("A2").Value = "20"

# number_of_sample: 2 -> 4
$statWs.Range("B2").NumberFormat = "@"
$statWs.Range("B2").Value = "4"

# number_of_cases: 1 -> 4
$statWs.Range("C2").NumberFormat = "@"
$statWs.Range("C2").Value = "4"

# number_of_study: 1 -> 2
$statWs.Range("D2").NumberFormat = "@"
$statWs.Range("D2").Value = "2"

# Drop the temporary Text number-format again so the cells end up with no
# explicit style applied (matches the original workbook's plain cells).
$statWs.Range("A2:D2").ClearFormats()

# --- StatOutput_Message sheet: the Cypher text logged for the stats query
#     is updated to reflect the new breed filter ---
$msgWs = $wb.Worksheets.Item("StatOutput_Message")
$newCypher = "MATCH (s:study) WITH COLLECT(DISTINCT(s.clinical_study_designation)) AS all_studies MATCH (d:demographic) WITH COLLECT(DISTINCT(d.breed)) AS all_breeds, COLLECT(DISTINCT(d.sex)) AS all_sexes, all_studies MATCH (d:diagnosis) WITH COLLECT(DISTINCT(d.disease_term)) AS all_diseases, all_breeds, all_sexes, all_studies MATCH (p:program)<-[*]-(s:study)<-[*]-(c:case)<--(demo:demographic), (c)<--(diag:diagnosis) WHERE demo.breed IN['German Shepherd Dog']  OPTIONAL MATCH (f:file)-[*]->(c), (samp:sample)-[*]->(c) WITH DISTINCT c AS c, p, s, demo, diag, f, samp RETURN count(DISTINCT(f)) as number_of_files , count(DISTINCT(samp)) as number_of_sample , count(DISTINCT(c.case_id)) as number_of_cases , count(DISTINCT(s.clinical_study_designation)) as number_of_study"
$msgWs.Range("A18").Value = $newCypher
